$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4052
$ws1.Range("F5").Value = 322
$ws1.Range("F8").Value = 31
$ws1.Range("F10").Value = 119
$ws1.Range("F11").Value = 294
$ws1.Range("F12").Value = 226
$ws1.Range("F13").Value = 2865
$ws1.Range("F14").Value = 127
$ws1.Range("F15").Value = 1278

# Sheet "全部类型" (All types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4052
$ws4.Range("F5").Value = 322
$ws4.Range("F9").Value = 31
$ws4.Range("F11").Value = 119
$ws4.Range("F12").Value = 294
$ws4.Range("F13").Value = 226
$ws4.Range("F14").Value = 2865
$ws4.Range("F15").Value = 127
$ws4.Range("F16").Value = 1278
